# Apply updated cryptocurrency price/volume data to the active worksheet.
# Values are written as plain text (via NumberFormat + Value) so that
# numeric-looking strings such as "246.02" or "35.239.23" remain text,
# matching the original inline-string cell type used for columns B-E.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = [ordered]@{
    "D2" = '35.239.23'
    "E2" = '  -0.80%  '
    "D3" = '1.902.94'
    "E3" = '  +1.23%  '
    "B5" = 'XRP'
    "C5" = 'https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp'
    "D5" = '0.695'
    "E5" = '  +9.94%  '
    "B6" = 'BNB'
    "C6" = 'https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb'
    "D6" = '246.02'
    "E6" = '  +1.19%  '
    "E7" = '  -0.40%  '
    "D8" = '42.09'
    "E8" = '  -2.45%  '
    "E9" = '  +5.11%  '
    "D10" = '53.19'
    "E10" = '  +11.78%  '
    "D11" = '0.0727'
    "E11" = '  +3.28%  '
    "D12" = '0.0996'
    "E12" = '  +0.28%  '
    "D13" = '2.179.84'
    "E13" = '  +1.36%  '
    "D14" = '12.32'
    "E14" = '  +1.14%  '
    "D15" = '0.710'
    "E15" = '  +3.46%  '
    "D16" = '1.897.23'
    "E16" = '  +0.66%  '
    "D17" = '4.84'
    "E17" = '  +1.25%  '
    "D18" = '35.255.81'
    "E18" = '  -0.74%  '
    "D19" = '72.47'
    "E19" = '  +1.63%  '
    "D20" = '0.0₃0822'
    "E20" = '  +1.89%  '
    "D21" = '241.13'
    "E21" = '  -1.06%  '
    "D22" = '12.60'
    "E22" = '  +1.56%  '
    "D23" = '4.84'
    "E23" = '  +0.06%  '
    "E24" = '  -0.39%  '
    "E25" = '  +2.22%  '
    "D26" = '2.33'
    "E26" = '  +15.14%  '
    "D27" = '169.88'
    "E27" = '  -0.80%  '
    "E28" = '  +3.32%  '
    "E29" = '  +4.93%  '
    "D30" = '18.41'
    "E30" = '  +3.16%  '
    "B32" = 'Filecoin'
    "C32" = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
    "D32" = '4.18'
    "E32" = '  +2.50%  '
    "B33" = 'ImmutableX'
    "C33" = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
    "D33" = '0.965'
    "E33" = '  -0.36%  '
    "D34" = '0.0572'
    "E34" = '  +1.33%  '
    "E35" = '  -0.47%  '
    "E36" = '  +0.70%  '
    "E37" = '  -0.57%  '
    "E38" = '  -1.30%  '
    "E39" = '  -1.27%  '
    "D40" = '0.0678'
    "E40" = '  +14.55%  '
    "E41" = '  +0.30%  '
    "E42" = '  +2.77%  '
    "D43" = '16.23'
    "E43" = '  +6.91%  '
    "D44" = '90.48'
    "E44" = '  +0.34%  '
    "D45" = '1.343.68'
    "D46" = '2.43'
    "E46" = '  +4.03%  '
    "D47" = '47.02'
    "E47" = '  +4.49%  '
    "D48" = '12.61'
    "E48" = '  -0.23%  '
    "D49" = '2.42'
    "E49" = '  -0.35%  '
    "D50" = '2.80'
    "E50" = '  +1.60%  '
    "E51" = '  -2.36%  '
}

foreach ($addr in $updates.Keys) {
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $updates[$addr]
}
